# Generate Report for Handoff
# Inserts a new row (for file 55d0b76c-f4eb-42dc-9129-283d748e1e3f) as the
# second-to-last data row on each of the three sheets (Overview, zh-cn,
# de-de), pushing the previous last row (b1d76ec3-...) down by one.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet 1: Overview  (columns A:D -> File Name, zh-cn, de-de, Latest Handoff Date)
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Overview")

$ws1.Rows("7:7").Insert()

$ws1.Range("A7").Value = "55d0b76c-f4eb-42dc-9129-283d748e1e3f.md"
$ws1.Range("B7").Value = "Ready for handoff"
$ws1.Range("C7").Value = "Ready for handoff"
$ws1.Range("D7").Value = "2016-03-25 09:42:47"

# Rebuild every hyperlink on this sheet (row insert does not auto-shift them)
$ws1.Range("A2").Hyperlinks.Delete()
$ws1.Hyperlinks.Add($ws1.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/0b436017fae00705ebfdff30373cfe324458c4fe/e2e/23c30e2f-eb80-43a7-bff5-73ee950a18ed.md", "", "", "23c30e2f-eb80-43a7-bff5-73ee950a18ed.md") | Out-Null
$ws1.Hyperlinks.Add($ws1.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/aa3a6d724c73717e7ccf08de513b9ab3a0181332/e2e/45094b65-7ed6-4509-89a3-262a321170a9.md", "", "", "45094b65-7ed6-4509-89a3-262a321170a9.md") | Out-Null
$ws1.Hyperlinks.Add($ws1.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/aa3a6d724c73717e7ccf08de513b9ab3a0181332/e2e/4ee34909-2f92-40b9-af95-7432e3091794.md", "", "", "4ee34909-2f92-40b9-af95-7432e3091794.md") | Out-Null
$ws1.Hyperlinks.Add($ws1.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/6aa43d4f27c5bc787eae7ed6a4380b2ecdb14139/e2e/c142d7d7-760c-43f0-9cfb-e7ca64b4e3b2.md", "", "", "c142d7d7-760c-43f0-9cfb-e7ca64b4e3b2.md") | Out-Null
$ws1.Hyperlinks.Add($ws1.Range("A6"), "https://github.com/OpenLocalizationTest/oltest/blob/850a98e9f44a71a081685c563f2c6032f17e820d/e2e/28d1a69d-8062-450e-a3a4-9ef63aa3000e.md", "", "", "28d1a69d-8062-450e-a3a4-9ef63aa3000e.md") | Out-Null
$ws1.Hyperlinks.Add($ws1.Range("A7"), "https://github.com/OpenLocalizationTest/oltest/blob/7ddcb9b484123d4303ffa78f9dd301c9e46ee1ad/e2e/55d0b76c-f4eb-42dc-9129-283d748e1e3f.md", "", "", "55d0b76c-f4eb-42dc-9129-283d748e1e3f.md") | Out-Null
$ws1.Hyperlinks.Add($ws1.Range("A8"), "https://github.com/OpenLocalizationTest/oltest/blob/d066a162032495ca6ccc28d39413caa152d8ea26/e2e/b1d76ec3-3e9b-4226-ab79-56cc1be0a550.md", "", "", "b1d76ec3-3e9b-4226-ab79-56cc1be0a550.md") | Out-Null

# ---------------------------------------------------------------------
# Sheet 2: zh-cn
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("zh-cn")

$ws2.Rows("7:7").Insert()

$ws2.Range("A7").Value = "55d0b76c-f4eb-42dc-9129-283d748e1e3f.md"
$ws2.Range("B7").Value = ".md"
$ws2.Range("C7").Value = "Ready for handoff"
$ws2.Range("D7").Value = "55d0b76c-f4eb-42dc-9129-283d748e1e3f.8c57edb29cfecc372566b892e601e2546d6cc719.zh-cn.xlf"
$ws2.Range("E7").Value = "2016-03-25 09:42:38"
$ws2.Range("H7").Value = "0001-01-01 00:00:00"
$ws2.Range("J7").Value = "Include"

$ws2.Range("A2").Hyperlinks.Delete()
$ws2.Hyperlinks.Add($ws2.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/0b436017fae00705ebfdff30373cfe324458c4fe/e2e/23c30e2f-eb80-43a7-bff5-73ee950a18ed.md", "", "", "23c30e2f-eb80-43a7-bff5-73ee950a18ed.md") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("D2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/3b72826d7f986d75c8cc9024ac010b52d240eddf/ol-handoff/OpenLocalizationTestOrg/oltest-zhcn-fly/yuwzho/ht/23c30e2f-eb80-43a7-bff5-73ee950a18ed.59325c62c2ee183de4bfb3d00499b6c0be81cf82.zh-cn.xlf", "", "", "23c30e2f-eb80-43a7-bff5-73ee950a18ed.59325c62c2ee183de4bfb3d00499b6c0be81cf82.zh-cn.xlf") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/aa3a6d724c73717e7ccf08de513b9ab3a0181332/e2e/45094b65-7ed6-4509-89a3-262a321170a9.md", "", "", "45094b65-7ed6-4509-89a3-262a321170a9.md") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("D3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/63859ceed7aa5a9106c939720e094ac5c4cb083b/ol-handoff/OpenLocalizationTestOrg/oltest-zhcn-fly/yuwzho/ht/45094b65-7ed6-4509-89a3-262a321170a9.5e4090b7461d36ef5f7142747c922b0b6c5c5709.zh-cn.xlf", "", "", "45094b65-7ed6-4509-89a3-262a321170a9.5e4090b7461d36ef5f7142747c922b0b6c5c5709.zh-cn.xlf") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/aa3a6d724c73717e7ccf08de513b9ab3a0181332/e2e/4ee34909-2f92-40b9-af95-7432e3091794.md", "", "", "4ee34909-2f92-40b9-af95-7432e3091794.md") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("D4"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/63859ceed7aa5a9106c939720e094ac5c4cb083b/ol-handoff/OpenLocalizationTestOrg/oltest-zhcn-fly/yuwzho/ht/4ee34909-2f92-40b9-af95-7432e3091794.63d9a3b2d2936c58ad261e2d1762983b6c608073.zh-cn.xlf", "", "", "4ee34909-2f92-40b9-af95-7432e3091794.63d9a3b2d2936c58ad261e2d1762983b6c608073.zh-cn.xlf") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/6aa43d4f27c5bc787eae7ed6a4380b2ecdb14139/e2e/c142d7d7-760c-43f0-9cfb-e7ca64b4e3b2.md", "", "", "c142d7d7-760c-43f0-9cfb-e7ca64b4e3b2.md") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("D5"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/f23d693728110ce20ce9accbe741245f12005295/ol-handoff/OpenLocalizationTestOrg/oltest-zhcn-fly/yuwzho/ht/c142d7d7-760c-43f0-9cfb-e7ca64b4e3b2.53239a2c51cb46048413ba0b433ec314bc9f21ae.zh-cn.xlf", "", "", "c142d7d7-760c-43f0-9cfb-e7ca64b4e3b2.53239a2c51cb46048413ba0b433ec314bc9f21ae.zh-cn.xlf") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("F5"), "https://github.com/OpenLocalizationTestOrg/oltest-zhcn-fly/blob/94cbca976ebd3815be978ee1a970d5e491bf125e/e2e/c142d7d7-760c-43f0-9cfb-e7ca64b4e3b2.md", "", "", "c142d7d7-760c-43f0-9cfb-e7ca64b4e3b2.md") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("G5"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/a9bbd621b6576eed0de886c03b196796d13eedca/ol-handback/OpenLocalizationTestOrg/oltest-zhcn-fly/yuwzho/ht/c142d7d7-760c-43f0-9cfb-e7ca64b4e3b2.53239a2c51cb46048413ba0b433ec314bc9f21ae.zh-cn.xlf", "", "", "c142d7d7-760c-43f0-9cfb-e7ca64b4e3b2.53239a2c51cb46048413ba0b433ec314bc9f21ae.zh-cn.xlf") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("A6"), "https://github.com/OpenLocalizationTest/oltest/blob/850a98e9f44a71a081685c563f2c6032f17e820d/e2e/28d1a69d-8062-450e-a3a4-9ef63aa3000e.md", "", "", "28d1a69d-8062-450e-a3a4-9ef63aa3000e.md") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("D6"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/09f9b2243eccb53f771c9a62c9bbe703f9e084f6/ol-handoff/OpenLocalizationTestOrg/oltest-zhcn-fly/yuwzho/ht/28d1a69d-8062-450e-a3a4-9ef63aa3000e.c20bb1386b847bb1bdbf93b611fbc6ac2eed13f3.zh-cn.xlf", "", "", "28d1a69d-8062-450e-a3a4-9ef63aa3000e.c20bb1386b847bb1bdbf93b611fbc6ac2eed13f3.zh-cn.xlf") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("A7"), "https://github.com/OpenLocalizationTest/oltest/blob/7ddcb9b484123d4303ffa78f9dd301c9e46ee1ad/e2e/55d0b76c-f4eb-42dc-9129-283d748e1e3f.md", "", "", "55d0b76c-f4eb-42dc-9129-283d748e1e3f.md") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("D7"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/d0c98b44574d8997d8a60f6a4ceaa3c8f18226a5/ol-handoff/OpenLocalizationTestOrg/oltest-zhcn-fly/yuwzho/ht/55d0b76c-f4eb-42dc-9129-283d748e1e3f.8c57edb29cfecc372566b892e601e2546d6cc719.zh-cn.xlf", "", "", "55d0b76c-f4eb-42dc-9129-283d748e1e3f.8c57edb29cfecc372566b892e601e2546d6cc719.zh-cn.xlf") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("A8"), "https://github.com/OpenLocalizationTest/oltest/blob/d066a162032495ca6ccc28d39413caa152d8ea26/e2e/b1d76ec3-3e9b-4226-ab79-56cc1be0a550.md", "", "", "b1d76ec3-3e9b-4226-ab79-56cc1be0a550.md") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("D8"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/4585a97f65489c29da60335cb320e9d4c3d0eb1f/ol-handoff/OpenLocalizationTestOrg/oltest-zhcn-fly/yuwzho/ht/b1d76ec3-3e9b-4226-ab79-56cc1be0a550.d6f4f91fa6d8dadb48aff48f0a1a5e3ea8762b04.zh-cn.xlf", "", "", "b1d76ec3-3e9b-4226-ab79-56cc1be0a550.d6f4f91fa6d8dadb48aff48f0a1a5e3ea8762b04.zh-cn.xlf") | Out-Null

# ---------------------------------------------------------------------
# Sheet 3: de-de
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("de-de")

$ws3.Rows("7:7").Insert()

$ws3.Range("A7").Value = "55d0b76c-f4eb-42dc-9129-283d748e1e3f.md"
$ws3.Range("B7").Value = ".md"
$ws3.Range("C7").Value = "Ready for handoff"
$ws3.Range("D7").Value = "55d0b76c-f4eb-42dc-9129-283d748e1e3f.8c57edb29cfecc372566b892e601e2546d6cc719.de-de.xlf"
$ws3.Range("E7").Value = "Ready for handoff"
$ws3.Range("H7").Value = "0001-01-01 00:00:00"
$ws3.Range("J7").Value = "Include"

$ws3.Range("A2").Hyperlinks.Delete()
$ws3.Hyperlinks.Add($ws3.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/0b436017fae00705ebfdff30373cfe324458c4fe/e2e/23c30e2f-eb80-43a7-bff5-73ee950a18ed.md", "", "", "23c30e2f-eb80-43a7-bff5-73ee950a18ed.md") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("D2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/e73ac60c7e95c8321c2d1b7bdeaade5ad5d206a3/ol-handoff/OpenLocalizationTestOrg/oltest-dede-fly/yuwzho/ht/23c30e2f-eb80-43a7-bff5-73ee950a18ed.59325c62c2ee183de4bfb3d00499b6c0be81cf82.de-de.xlf", "", "", "23c30e2f-eb80-43a7-bff5-73ee950a18ed.59325c62c2ee183de4bfb3d00499b6c0be81cf82.de-de.xlf") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/aa3a6d724c73717e7ccf08de513b9ab3a0181332/e2e/45094b65-7ed6-4509-89a3-262a321170a9.md", "", "", "45094b65-7ed6-4509-89a3-262a321170a9.md") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("D3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/bf8671ecd7a2860b25318d0c9a5e46db46acbb67/ol-handoff/OpenLocalizationTestOrg/oltest-dede-fly/yuwzho/ht/45094b65-7ed6-4509-89a3-262a321170a9.5e4090b7461d36ef5f7142747c922b0b6c5c5709.de-de.xlf", "", "", "45094b65-7ed6-4509-89a3-262a321170a9.5e4090b7461d36ef5f7142747c922b0b6c5c5709.de-de.xlf") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/aa3a6d724c73717e7ccf08de513b9ab3a0181332/e2e/4ee34909-2f92-40b9-af95-7432e3091794.md", "", "", "4ee34909-2f92-40b9-af95-7432e3091794.md") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("D4"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/bf8671ecd7a2860b25318d0c9a5e46db46acbb67/ol-handoff/OpenLocalizationTestOrg/oltest-dede-fly/yuwzho/ht/4ee34909-2f92-40b9-af95-7432e3091794.63d9a3b2d2936c58ad261e2d1762983b6c608073.de-de.xlf", "", "", "4ee34909-2f92-40b9-af95-7432e3091794.63d9a3b2d2936c58ad261e2d1762983b6c608073.de-de.xlf") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/6aa43d4f27c5bc787eae7ed6a4380b2ecdb14139/e2e/c142d7d7-760c-43f0-9cfb-e7ca64b4e3b2.md", "", "", "c142d7d7-760c-43f0-9cfb-e7ca64b4e3b2.md") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("D5"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/2f45fc3ab4e23d8b432218a56930b80d1059565e/ol-handoff/OpenLocalizationTestOrg/oltest-dede-fly/yuwzho/ht/c142d7d7-760c-43f0-9cfb-e7ca64b4e3b2.53239a2c51cb46048413ba0b433ec314bc9f21ae.de-de.xlf", "", "", "c142d7d7-760c-43f0-9cfb-e7ca64b4e3b2.53239a2c51cb46048413ba0b433ec314bc9f21ae.de-de.xlf") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("F5"), "https://github.com/OpenLocalizationTestOrg/oltest-dede-fly/blob/d41f3d4d2c05f27b6e6f54ffde042cf1711f4e3e/e2e/c142d7d7-760c-43f0-9cfb-e7ca64b4e3b2.md", "", "", "c142d7d7-760c-43f0-9cfb-e7ca64b4e3b2.md") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("G5"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/10ec6855f31726a56b32c9066ff87505a9b328c2/ol-handback/OpenLocalizationTestOrg/oltest-dede-fly/yuwzho/ht/c142d7d7-760c-43f0-9cfb-e7ca64b4e3b2.53239a2c51cb46048413ba0b433ec314bc9f21ae.de-de.xlf", "", "", "c142d7d7-760c-43f0-9cfb-e7ca64b4e3b2.53239a2c51cb46048413ba0b433ec314bc9f21ae.de-de.xlf") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("A6"), "https://github.com/OpenLocalizationTest/oltest/blob/850a98e9f44a71a081685c563f2c6032f17e820d/e2e/28d1a69d-8062-450e-a3a4-9ef63aa3000e.md", "", "", "28d1a69d-8062-450e-a3a4-9ef63aa3000e.md") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("D6"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/a732c26360a4d20a833085cb3a9187250bc94029/ol-handoff/OpenLocalizationTestOrg/oltest-dede-fly/yuwzho/ht/28d1a69d-8062-450e-a3a4-9ef63aa3000e.c20bb1386b847bb1bdbf93b611fbc6ac2eed13f3.de-de.xlf", "", "", "28d1a69d-8062-450e-a3a4-9ef63aa3000e.c20bb1386b847bb1bdbf93b611fbc6ac2eed13f3.de-de.xlf") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("A7"), "https://github.com/OpenLocalizationTest/oltest/blob/7ddcb9b484123d4303ffa78f9dd301c9e46ee1ad/e2e/55d0b76c-f4eb-42dc-9129-283d748e1e3f.md", "", "", "55d0b76c-f4eb-42dc-9129-283d748e1e3f.md") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("D7"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/5443efcd81ae89ea81f5f8e81005846c16abd814/ol-handoff/OpenLocalizationTestOrg/oltest-dede-fly/yuwzho/ht/55d0b76c-f4eb-42dc-9129-283d748e1e3f.8c57edb29cfecc372566b892e601e2546d6cc719.de-de.xlf", "", "", "55d0b76c-f4eb-42dc-9129-283d748e1e3f.8c57edb29cfecc372566b892e601e2546d6cc719.de-de.xlf") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("A8"), "https://github.com/OpenLocalizationTest/oltest/blob/d066a162032495ca6ccc28d39413caa152d8ea26/e2e/b1d76ec3-3e9b-4226-ab79-56cc1be0a550.md", "", "", "b1d76ec3-3e9b-4226-ab79-56cc1be0a550.md") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("D8"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/8f629591920758487f3ceaedb4cf2d2957ca6172/ol-handoff/OpenLocalizationTestOrg/oltest-dede-fly/yuwzho/ht/b1d76ec3-3e9b-4226-ab79-56cc1be0a550.d6f4f91fa6d8dadb48aff48f0a1a5e3ea8762b04.de-de.xlf", "", "", "b1d76ec3-3e9b-4226-ab79-56cc1be0a550.d6f4f91fa6d8dadb48aff48f0a1a5e3ea8762b04.de-de.xlf") | Out-Null

Write-Output "Report generated for handoff."
